$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing summary formula in G379 (sum of F129 + F130)
$ws.Range("G379").Formula = "=F129+F130"

# Re-apply the AutoFilter: narrow CENTROS (col B) to DOÑA CARMEN / SARMIENTO,
# and change the indicator (col C) filter to "9.1-Cobertura atención salud mental".
# This also recomputes which data rows are hidden vs visible.
$ws.Range("A1:L375").AutoFilter(2, @("DOÑA CARMEN","SARMIENTO"), 7)
$ws.Range("A1:L375").AutoFilter(3, @("9.1-Cobertura atención salud mental"), 7)

# Update the active selection to reflect the newly visible rows.
$ws.Range("F197:F198").Select()

Write-Host "edit complete"
